$wb = $excel.ActiveWorkbook

# ----- AR sheet -----
$wsAR = $wb.Worksheets.Item("AR")
$wsAR.Range("B2").Value = 0.006270035935739697
$wsAR.Range("B3").Value = 0.7913659530258792
$wsAR.Range("B4").Value = 0.09797099665097854
$wsAR.Range("B5").Value = "[0.9999999999999999, 0.027519311876138233, -0.0510786890162894, -0.09506615120389852, -0.1359016308868483, 0.0748412585706099, 0.029298285440112162, 0.05585394029129564, 0.025298370041328528, -0.0300446661980944, -0.028319987521129998, -0.06437384264461887, 0.003492010995920257, 0.04470859488424577, 0.06553932197021563, 0.029596804916165165, -0.02486665766923546, -0.050434194654770365, -0.04090466237308027, 0.011400893517586221]"

# ----- SETAR sheet -----
$wsSETAR = $wb.Worksheets.Item("SETAR")
$wsSETAR.Range("B4").Value = -0.1736855735702523
$wsSETAR.Range("B5").Value = 0.5771728172191997
$wsSETAR.Range("B6").Value = 0.05726102397873654
$wsSETAR.Range("B7").Value = 0.1630749003618141
$wsSETAR.Range("B8").Value = 0.6136283770559584
$wsSETAR.Range("B9").Value = 0.05754291918216003
$wsSETAR.Range("B10").Value = "[1.0, 0.05285221498209817, -0.01420217580304059, -0.014710379035803932, -0.02143248630957043, 0.01643521218672918, -0.010207665063015047, 0.0213482407275454, -0.00856893386069431, -0.030367691149575625, -0.004297175689131309, -0.021358863509527592, -0.016861521769549012, -0.000865443254141501, 0.04458572686571549, 0.00786045057787506, -0.005546749481335231, 0.022367006114971232, -0.015106608755289583, 0.007787175937302521]"

# ----- GARCH sheet -----
$wsGARCH = $wb.Worksheets.Item("GARCH")
$wsGARCH.Range("B2").Value = 0.0001306639538657085
$wsGARCH.Range("B3").Value = 0.000000001094242180379503
$wsGARCH.Range("B4").Value = 0.0008793495878507832
$wsGARCH.Range("B5").Value = 0.9990651578338104
$wsGARCH.Range("B6").Value = "[1.0, 0.028540450934653808, -0.04253255726384514, -0.08768965618830829, -0.13996844894994742, 0.09175255283912845, 0.02014256668219339, 0.06757728614698867, 0.03202650222886918, -0.031229687706606845, -0.04134583334353597, -0.06302856623731395, 0.008599507863616752, 0.03490754169297713, 0.06963961681611058, 0.02801339035048576, -0.03609694384209265, -0.04369170970762815, -0.05273206232509614, 0.0026995178662624474]"

# ----- TARCH sheet -----
$wsTARCH = $wb.Worksheets.Item("TARCH")
$wsTARCH.Range("B2").Value = -0.0007090862305013981
$wsTARCH.Range("B3").Value = 0.1045784134493644
$wsTARCH.Range("B4").Value = 0.02002084811503992
$wsTARCH.Range("B5").Value = 0
$wsTARCH.Range("B6").Value = "[1.0, -0.003927731358294636, -0.03972892414594325, -0.0808667625881799, -0.13698441887718352, 0.09609221935364978, 0.016421197551761052, 0.06777860249973515, 0.03245472545829643, -0.028019236275828325, -0.035634090000457364, -0.05817488322454031, 0.01189164111954234, 0.03506729135024723, 0.06812148640902077, 0.029768608518841323, -0.03292478901786283, -0.03782260469969652, -0.048052073335654075, 0.005582789216852781]"
$wsTARCH.Range("B7").Value = 0.04748264745024922

# ----- AR_TARCH sheet -----
$wsARTARCH = $wb.Worksheets.Item("AR_TARCH")
$wsARTARCH.Range("B2").Value = 0.006370006642671128
$wsARTARCH.Range("B3").Value = 0.09477549752665887
$wsARTARCH.Range("B4").Value = 0.0112803772086497
$wsARTARCH.Range("B5").Value = 0
$wsARTARCH.Range("B6").Value = "[1.0, -0.0007632531981645069, -0.050040660163139544, -0.09069172429952296, -0.13607836214244434, 0.07929928757085285, 0.025879841130390012, 0.05496972443405766, 0.025095690736669594, -0.030799100600430236, -0.025109748421111636, -0.06319761703726988, 0.0025698986699293066, 0.04379026239652982, 0.0627432121610422, 0.02924235230031501, -0.02342059656590508, -0.04809899344282776, -0.039597172199580234, 0.010276307471163648]"
$wsARTARCH.Range("B7").Value = 0.04179881657550608
$wsARTARCH.Range("B9").Value = 0.7923285543529349
